$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new fertilizer type entry as the next row in the list
$ws.Range("A12").Value = "10-34-0"

# Move selection to the next empty cell, mirroring Excel's post-entry behavior
$ws.Range("A13").Select()
